$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old trailing rows (18-20) first so row 17 stays put while we rewrite it
$ws.Rows("18:20").Delete()

# New header for column D
$ws.Range("D1").Value = "ITI"

# Column C (ConditionType) updates + new column D (ITI) values, rows 2-17
$ws.Range("C2").Value = 4
$ws.Range("D2").Value = 6

$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 8

$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 6

$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 6

$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 7

$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 6

$ws.Range("C8").Value = 2
$ws.Range("D8").Value = 8

$ws.Range("C9").Value = 1
$ws.Range("D9").Value = 8

$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 8

$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 7

$ws.Range("C12").Value = 4
$ws.Range("D12").Value = 8

$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 6

$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 8

$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 7

$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 6

$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 7

# Match the saved selection state from the diff
$ws.Range("I7").Select() | Out-Null
